$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "28.896.53"
Set-TextValue $ws.Range("E2") "  -1.73%  "
Set-TextValue $ws.Range("D3") "1.833.76"
Set-TextValue $ws.Range("E3") "  -1.96%  "
Set-TextValue $ws.Range("D4") "0.9993"
Set-TextValue $ws.Range("E4") "  -0.11%  "
Set-TextValue $ws.Range("D5") "245.14"
Set-TextValue $ws.Range("E5") "  +0.64%  "
Set-TextValue $ws.Range("D6") "0.6884"
Set-TextValue $ws.Range("D7") "0.9998"
Set-TextValue $ws.Range("D8") "0.07691"
Set-TextValue $ws.Range("E8") "  -2.91%  "
Set-TextValue $ws.Range("D9") "0.3054"
Set-TextValue $ws.Range("E9") "  -2.52%  "
Set-TextValue $ws.Range("D10") "23.42"
Set-TextValue $ws.Range("E10") "  -4.52%  "
Set-TextValue $ws.Range("D11") "0.07819"
Set-TextValue $ws.Range("E11") "  +0.08%  "
Set-TextValue $ws.Range("D12") "1.831.35"
Set-TextValue $ws.Range("E12") "  -3.31%  "
Set-TextValue $ws.Range("D13") "5.098"
Set-TextValue $ws.Range("E13") "  -1.44%  "
Set-TextValue $ws.Range("D14") "90.44"
Set-TextValue $ws.Range("E14") "  -3.46%  "
Set-TextValue $ws.Range("D15") "0.6822"
Set-TextValue $ws.Range("E15") "  -2.95%  "
Set-TextValue $ws.Range("D16") "6.444"
Set-TextValue $ws.Range("E16") "  -1.27%  "
Set-TextValue $ws.Range("D17") "0.000008306"
Set-TextValue $ws.Range("E17") "  -1.45%  "
Set-TextValue $ws.Range("D18") "28.884.24"
Set-TextValue $ws.Range("E18") "  -1.98%  "
Set-TextValue $ws.Range("D19") "242.64"
Set-TextValue $ws.Range("E19") "  -3.79%  "
Set-TextValue $ws.Range("D20") "2.074.99"
Set-TextValue $ws.Range("E20") "  -3.53%  "
Set-TextValue $ws.Range("D21") "12.74"
Set-TextValue $ws.Range("E21") "  -2.92%  "
Set-TextValue $ws.Range("D22") "0.9999"
Set-TextValue $ws.Range("E22") "  -0.01%  "
Set-TextValue $ws.Range("D23") "7.466"
Set-TextValue $ws.Range("D24") "0.9996"
Set-TextValue $ws.Range("E24") "  -0.09%  "
Set-TextValue $ws.Range("E25") "  -4.93%  "
Set-TextValue $ws.Range("D26") "161.35"
Set-TextValue $ws.Range("E26") "  -0.17%  "
Set-TextValue $ws.Range("D27") "8.816"
Set-TextValue $ws.Range("E27") "  -2.20%  "
Set-TextValue $ws.Range("D28") "18.20"
Set-TextValue $ws.Range("E28") "  -3.27%  "
Set-TextValue $ws.Range("D29") "1.546"
Set-TextValue $ws.Range("E29") "  +2.85%  "
Set-TextValue $ws.Range("D30") "4.216"
Set-TextValue $ws.Range("E30") "  -2.23%  "
Set-TextValue $ws.Range("D31") "4.156"
Set-TextValue $ws.Range("E31") "  -2.42%  "
Set-TextValue $ws.Range("D32") "1.184"
Set-TextValue $ws.Range("E32") "  -2.55%  "
Set-TextValue $ws.Range("D33") "0.05113"
Set-TextValue $ws.Range("E33") "  -2.88%  "
Set-TextValue $ws.Range("D34") "0.7658"
Set-TextValue $ws.Range("E34") "  +1.67%  "
Set-TextValue $ws.Range("E35") "  -2.94%  "
Set-TextValue $ws.Range("E36") "  -3.42%  "
Set-TextValue $ws.Range("E37") "  -0.80%  "
Set-TextValue $ws.Range("D38") "0.01849"
Set-TextValue $ws.Range("E38") "  -1.49%  "
Set-TextValue $ws.Range("D39") "1.224.84"
Set-TextValue $ws.Range("E39") "  -4.27%  "
Set-TextValue $ws.Range("E40") "  -2.50%  "
Set-TextValue $ws.Range("D41") "0.9413"
Set-TextValue $ws.Range("D42") "108.71"
Set-TextValue $ws.Range("E42") "  -0.44%  "
Set-TextValue $ws.Range("D43") "0.9994"
Set-TextValue $ws.Range("D44") "5.710"
Set-TextValue $ws.Range("E44") "  -5.22%  "
Set-TextValue $ws.Range("E45") "  -3.12%  "
Set-TextValue $ws.Range("D46") "9.542"
Set-TextValue $ws.Range("E46") "  -0.67%  "
Set-TextValue $ws.Range("E47") "  -0.60%  "
Set-TextValue $ws.Range("D48") "1.974.24"
Set-TextValue $ws.Range("E48") "  -3.48%  "
Set-TextValue $ws.Range("D49") "64.32"
Set-TextValue $ws.Range("E49") "  -9.32%  "
Set-TextValue $ws.Range("D50") "1.748"
Set-TextValue $ws.Range("E50") "  -3.12%  "
Set-TextValue $ws.Range("D51") "0.4189"
Set-TextValue $ws.Range("E51") "  -2.53%  "
